$wb = $excel.ActiveWorkbook

# --- Stages sheet: rename the placeholder "TBD" stage to "IRGen" ---
$wsStages = $wb.Worksheets.Item("Stages")
$wsStages.Range("A5").Value = "IRGen"
$wsStages.Range("A6").Select() | Out-Null

# --- Errors sheet: append two new IRGen errors to the Table1 table ---
$ws = $wb.Worksheets.Item("Errors")
$tbl = $ws.ListObjects.Item(1)

$row30 = $tbl.ListRows.Add()
$ws.Range("A30").Value = "Error"
$ws.Range("B30").Value = 1
$ws.Range("C30").Value = "IRGen"
$ws.Range("E30").Formula = "= (XLOOKUP(`$C30,Stages!`$A:`$A,Stages!`$B:`$B)+`$B30)"
$ws.Range("F30").Formula = "= LEFT(A30,1)&E30"

$row31 = $tbl.ListRows.Add()
$ws.Range("A31").Value = "Error"
$ws.Range("B31").Value = 2
$ws.Range("C31").Value = "IRGen"
$ws.Range("E31").Formula = "= (XLOOKUP(`$C31,Stages!`$A:`$A,Stages!`$B:`$B)+`$B31)"
$ws.Range("F31").Formula = "= LEFT(A31,1)&E31"

# Errors stays the active sheet with D31 selected, matching the post-edit state
$ws.Activate() | Out-Null
$ws.Range("D31").Select() | Out-Null
